$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the value to be written as literal text (so it becomes a shared
    # string instead of a numeric/date cell) without leaving any trace of a
    # formula or custom number format behind.
    $ws.Range($range).Formula = '="' + $text + '"'
    $ws.Range($range).Copy() | Out-Null
    $ws.Range($range).PasteSpecial(-4163) | Out-Null
}

# Row 2 - timestamps/ID refreshed (stays a "Completed" row)
Set-TextValue "A2" "04/23/2024 14:13:33"
Set-TextValue "B2" "04/23/2024 14:14:05"
Set-TextValue "C2" "183"
$ws.Range("D2").Value = "Completed"

# Row 3 - previously the failed run, now a new "Completed" run (Automation
# Exercise Logout). Drop the Error Message / Screenshot Path columns.
Set-TextValue "A3" "04/23/2024 14:14:05"
Set-TextValue "B3" "04/23/2024 14:14:35"
Set-TextValue "C3" "174"
$ws.Range("D3").Value = "Completed"
$ws.Range("E3:F3").ClearContents()

# Row 4 - previously a "Completed" run, now the failed run with its message
# and screenshot path.
Set-TextValue "A4" "04/23/2024 14:14:35"
Set-TextValue "B4" "04/23/2024 14:14:51"
Set-TextValue "C4" "155"
$ws.Range("D4").Value = "Failed"
$ws.Range("E4").Value = "Personal Info Generated did not Match"
$ws.Range("F4").Value = "C:\Users\pc\Desktop\screenshots\ExceptionScreenshot_240423.021451.png"
